# Updated league table for GW26.
# Enter this gameweek's scores for row 23 (week 25), following the same
# pattern used for the previous rows (20-22): copy the formatting from the
# row above, then fill in Eren/Mert/Arda's points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20:D20").Copy()
$ws.Range("B23:D23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B23").Value = 87
$ws.Range("C23").Value = 99
$ws.Range("D23").Value = 103

$ws.Range("C24").Select()
